# Change the Runmode for the "AuthoringTest" and "AuthoringAppreciateTest"
# test cases from "Y" to "N" so that only "AuthoringDeleteTest" and
# "AuthoringProfileCommentsTest" are run (see commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
